$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Selecting the inline picture first and renaming it through the live
# Selection object (rather than the InlineShapes wrapper returned directly
# off HeaderFooter.Range) is what makes the rename reliably reach the
# underlying picture for every story, footers included.
function Rename-InlinePicture($range, $newName) {
    $pic = $range.InlineShapes.Item(1)
    $pic.Select()
    $sel = $word.Selection
    $sel.InlineShapes.Item(1).Name = $newName
}

# --- Footers: Pearson logo, rename image1.png -> image2.png ---
# Sections(1).Footers(1) == word/footer2.xml ; Sections(1).Footers(2) == word/footer1.xml
Rename-InlinePicture $sec.Footers.Item(1).Range "image2.png"
Rename-InlinePicture $sec.Footers.Item(2).Range "image2.png"

# --- Headers: BTEC logo, rename image2.jpg -> image1.jpg ---
# Sections(1).Headers(1) == word/header2.xml ; Sections(1).Headers(2) == word/header1.xml
Rename-InlinePicture $sec.Headers.Item(1).Range "image1.jpg"
Rename-InlinePicture $sec.Headers.Item(2).Range "image1.jpg"

Write-Host "Renamed footer and header inline picture names"
